$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is a numeric-looking string (e.g. "595.18") need to
# be pre-formatted as Text; otherwise Excel auto-converts the assigned string into a
# real number (losing the exact decimal text and changing the stored cell type).
$textCells = $excel.Union($ws.Range("D4"), $ws.Range("D5"), $ws.Range("D6"), $ws.Range("D8"), $ws.Range("D10"), $ws.Range("D12"), $ws.Range("D13"), $ws.Range("D20"), $ws.Range("D21"), $ws.Range("D23"), $ws.Range("D25"), $ws.Range("D27"), $ws.Range("D30"), $ws.Range("D32"), $ws.Range("D33"), $ws.Range("D34"), $ws.Range("D36"), $ws.Range("D37"), $ws.Range("D38"), $ws.Range("D40"), $ws.Range("D41"), $ws.Range("D42"), $ws.Range("D43"), $ws.Range("D44"), $ws.Range("D47"), $ws.Range("D48"), $ws.Range("D49"), $ws.Range("D51"))
$textCells.NumberFormat = "@"

$ws.Range("D2").Value = "64.340.19"
$ws.Range("E2").Value = "  -1.65%  "

$ws.Range("D3").Value = "3.127.08"
$ws.Range("E3").Value = "  -2.30%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").Value = "595.18"
$ws.Range("E5").Value = "  -0.68%  "

$ws.Range("D6").Value = "159.01"
$ws.Range("E6").Value = "  +3.34%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "0.543"
$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").Value = "3.125.69"
$ws.Range("E9").Value = "  -2.33%  "

$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  -5.10%  "

$ws.Range("E11").Value = "  -2.95%  "

$ws.Range("D12").Value = "0.455"
$ws.Range("E12").Value = "  -3.83%  "

$ws.Range("D13").Value = "37.44"
$ws.Range("E13").Value = "  -4.68%  "

$ws.Range("E14").Value = "  -5.92%  "

$ws.Range("D15").Value = "3.640.58"
$ws.Range("E15").Value = "  -2.45%  "

$ws.Range("E16").Value = "  -1.42%  "

$ws.Range("E17").Value = "  -1.91%  "

$ws.Range("D18").Value = "64.211.43"
$ws.Range("E18").Value = "  -1.36%  "

$ws.Range("D19").Value = "3.120.42"
$ws.Range("E19").Value = "  -2.58%  "

$ws.Range("D20").Value = "479.89"
$ws.Range("E20").Value = "  -1.12%  "

$ws.Range("D21").Value = "14.61"
$ws.Range("E21").Value = "  -3.14%  "

$ws.Range("E22").Value = "  -7.14%  "

$ws.Range("D23").Value = "7.64"
$ws.Range("E23").Value = "  -3.59%  "

$ws.Range("E24").Value = "  +1.40%  "

$ws.Range("D25").Value = "13.07"
$ws.Range("E25").Value = "  -6.04%  "

$ws.Range("D27").Value = "10.59"
$ws.Range("E27").Value = "  +7.66%  "

$ws.Range("E28").Value = "  -0.24%  "

$ws.Range("E29").Value = "  +1.27%  "

$ws.Range("D30").Value = "2.71"
$ws.Range("E30").Value = "  -2.82%  "

$ws.Range("E31").Value = "  -2.90%  "

$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  -0.35%  "

$ws.Range("D33").Value = "0.114"
$ws.Range("E33").Value = "  -6.10%  "

$ws.Range("D34").Value = "27.47"
$ws.Range("E34").Value = "  -4.03%  "

$ws.Range("D35").Value = "0.0" + [char]0x2083 + "0853"
$ws.Range("E35").Value = "  -5.39%  "

$ws.Range("D36").Value = "1.07"
$ws.Range("E36").Value = "  -2.53%  "

$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "3.33"
$ws.Range("E37").Value = "  -8.46%  "

$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").Value = "6.07"
$ws.Range("E38").Value = "  -4.61%  "

$ws.Range("E39").Value = "  -5.27%  "

$ws.Range("D40").Value = "51.19"
$ws.Range("E40").Value = "  -1.54%  "

$ws.Range("D41").Value = "454.11"
$ws.Range("E41").Value = "  -5.27%  "

$ws.Range("D42").Value = "9.19"
$ws.Range("E42").Value = "  -2.66%  "

$ws.Range("D43").Value = "0.294"
$ws.Range("E43").Value = "  -2.71%  "

$ws.Range("D44").Value = "0.0368"
$ws.Range("E44").Value = "  -4.06%  "

$ws.Range("E45").Value = "  -0.16%  "

$ws.Range("D46").Value = "2.847.93"
$ws.Range("E46").Value = "  -3.56%  "

$ws.Range("D47").Value = "40.11"
$ws.Range("E47").Value = "  +3.38%  "

$ws.Range("D48").Value = "130.83"
$ws.Range("E48").Value = "  -0.57%  "

$ws.Range("D49").Value = "25.90"
$ws.Range("E49").Value = "  +1.06%  "

$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("D51").Value = "2.28"
$ws.Range("E51").Value = "  -2.34%  "

# Restore the default (General) style on those cells so only the cell's stored
# type/text changes - no stray number-format style is left applied to the cell.
$textCells.Style = "Normal"
